$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.701.32"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.671.71"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "599.70"
$ws.Range("E5").Value = "  -1.24%  "

# Row 6 - Solana
$ws.Range("D6").Value = "156.69"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +5.04%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +4.54%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.57%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  -2.22%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  -0.30%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "29.31"
$ws.Range("E13").Value = "  -3.32%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "0.0000197"
$ws.Range("E14").Value = "  -2.13%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.151.54"
$ws.Range("E15").Value = "  -0.83%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.543.24"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.679.96"
$ws.Range("E17").Value = "  -0.70%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "12.86"
$ws.Range("E18").Value = "  +1.70%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "4.79"

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  +0.20%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "352.16"
$ws.Range("E21").Value = "  -1.79%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.04%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "69.71"
$ws.Range("E23").Value = "  -1.26%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +5.31%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  -2.36%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "  -2.07%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -2.52%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  -5.36%  "

# Row 29 - Aptos
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  -4.63%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.05%  "

# Row 31 - was PancakeSwap, now Bittensor
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "528.73"
$ws.Range("E31").Value = "  -2.23%  "

# Row 32 - was Bittensor, now PancakeSwap
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").Value = "  -3.07%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").Value = "  -1.83%  "

# Row 34 - was RenderToken, now NEARProtocol
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "5.51"
$ws.Range("E34").Value = "  +1.97%  "

# Row 35 - was NEARProtocol, now RenderToken
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value = "6.47"
$ws.Range("E35").Value = "  -3.24%  "

# Row 36 - was PolygonEcosystemToken, now EthereumClassic
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "20.60"
$ws.Range("E36").Value = "  -1.28%  "

# Row 37 - was EthereumClassic, now PolygonEcosystemToken
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "0.422"
$ws.Range("E37").Value = "  -2.46%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.05%  "

# Row 39 - Monero
$ws.Range("D39").Value = "159.20"
$ws.Range("E39").Value = "  -2.74%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  -2.68%  "

# Row 42 - Aave
$ws.Range("D42").Value = "163.60"
$ws.Range("E42").Value = "  -4.28%  "

# Row 43 - Filecoin
$ws.Range("D43").Value = "4.13"
$ws.Range("E43").Value = "  -1.17%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "2.33"
$ws.Range("E44").Value = "  +2.64%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "0.0610"

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "22.75"
$ws.Range("E46").Value = "  -3.18%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "0.638"
$ws.Range("E47").Value = "  -2.53%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  -3.12%  "

# Row 49 - BabyDogeCoin
$ws.Range("D49").Value = "0.0₆0263"
$ws.Range("E49").Value = "  +14.33%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +0.72%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "20.16"
$ws.Range("E51").Value = "  -3.68%  "
